$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E51').NumberFormat = '@'

$ws.Range('D2').Value = '67.384.76'
$ws.Range('E2').Value = '  -3.07%  '
$ws.Range('D3').Value = '3.498.42'
$ws.Range('E3').Value = '  -4.60%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '601.86'
$ws.Range('E5').Value = '  -3.58%  '
$ws.Range('D6').Value = '149.43'
$ws.Range('E6').Value = '  -5.92%  '
$ws.Range('D7').Value = '3.497.94'
$ws.Range('E7').Value = '  -4.53%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.481'
$ws.Range('E9').Value = '  -3.24%  '
$ws.Range('E10').Value = '  -4.66%  '
$ws.Range('D11').Value = '6.94'
$ws.Range('E11').Value = '  -3.99%  '
$ws.Range('D12').Value = '0.422'
$ws.Range('E12').Value = '  -4.36%  '
$ws.Range('E13').Value = '  -5.68%  '
$ws.Range('D14').Value = '4.087.88'
$ws.Range('E14').Value = '  -4.64%  '
$ws.Range('D15').Value = '31.52'
$ws.Range('E15').Value = '  -2.46%  '
$ws.Range('D16').Value = '3.491.91'
$ws.Range('E16').Value = '  -5.29%  '
$ws.Range('D17').Value = '67.258.10'
$ws.Range('E17').Value = '  -3.28%  '
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').Value = '14.98'
$ws.Range('E20').Value = '  -5.74%  '
$ws.Range('D21').Value = '446.11'
$ws.Range('E21').Value = '  -4.83%  '
$ws.Range('E22').Value = '  -13.06%  '
$ws.Range('D23').Value = '0.618'
$ws.Range('E23').Value = '  -5.10%  '
$ws.Range('E24').Value = '  -2.93%  '
$ws.Range('E25').Value = '  +5.75%  '
$ws.Range('D27').Value = '3.635.76'
$ws.Range('E27').Value = '  -4.66%  '
$ws.Range('E28').Value = '  -9.36%  '
$ws.Range('D29').Value = '8.19'
$ws.Range('E29').Value = '  -5.23%  '
$ws.Range('E30').Value = '  -5.45%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  -7.40%  '
$ws.Range('E33').Value = '  +1.23%  '
$ws.Range('D34').Value = '25.67'
$ws.Range('E34').Value = '  -3.56%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '6.06'
$ws.Range('E35').Value = '  -5.13%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.486.77'
$ws.Range('E36').Value = '  -5.05%  '
$ws.Range('E37').Value = '  -6.60%  '
$ws.Range('D38').Value = '7.97'
$ws.Range('E38').Value = '  -3.67%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '2.19'
$ws.Range('E41').Value = '  -1.23%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = '174.41'
$ws.Range('E42').Value = '  -2.27%  '
$ws.Range('D43').Value = '0.0874'
$ws.Range('E43').Value = '  -2.09%  '
$ws.Range('D44').Value = '5.39'
$ws.Range('E44').Value = '  -7.04%  '
$ws.Range('D45').Value = '0.880'
$ws.Range('E45').Value = '  -4.66%  '
$ws.Range('D46').Value = '45.45'
$ws.Range('E46').Value = '  -3.09%  '
$ws.Range('E47').Value = '  +6.65%  '
$ws.Range('D48').Value = '27.23'
$ws.Range('E48').Value = '  -6.79%  '
$ws.Range('D49').Value = '2.55'
$ws.Range('E49').Value = '  -5.53%  '
$ws.Range('D50').Value = '7.52'
$ws.Range('E50').Value = '  -4.32%  '
$ws.Range('E51').Value = '  -3.95%  '

$ws.Range('D2').ClearFormats()
$ws.Range('E2').ClearFormats()
$ws.Range('D3').ClearFormats()
$ws.Range('E3').ClearFormats()
$ws.Range('E4').ClearFormats()
$ws.Range('D5').ClearFormats()
$ws.Range('E5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('E6').ClearFormats()
$ws.Range('D7').ClearFormats()
$ws.Range('E7').ClearFormats()
$ws.Range('E8').ClearFormats()
$ws.Range('D9').ClearFormats()
$ws.Range('E9').ClearFormats()
$ws.Range('E10').ClearFormats()
$ws.Range('D11').ClearFormats()
$ws.Range('E11').ClearFormats()
$ws.Range('D12').ClearFormats()
$ws.Range('E12').ClearFormats()
$ws.Range('E13').ClearFormats()
$ws.Range('D14').ClearFormats()
$ws.Range('E14').ClearFormats()
$ws.Range('D15').ClearFormats()
$ws.Range('E15').ClearFormats()
$ws.Range('D16').ClearFormats()
$ws.Range('E16').ClearFormats()
$ws.Range('D17').ClearFormats()
$ws.Range('E17').ClearFormats()
$ws.Range('E18').ClearFormats()
$ws.Range('E19').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('E20').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('E21').ClearFormats()
$ws.Range('E22').ClearFormats()
$ws.Range('D23').ClearFormats()
$ws.Range('E23').ClearFormats()
$ws.Range('E24').ClearFormats()
$ws.Range('E25').ClearFormats()
$ws.Range('D27').ClearFormats()
$ws.Range('E27').ClearFormats()
$ws.Range('E28').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('E29').ClearFormats()
$ws.Range('E30').ClearFormats()
$ws.Range('E31').ClearFormats()
$ws.Range('E32').ClearFormats()
$ws.Range('E33').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('E34').ClearFormats()
$ws.Range('D35').ClearFormats()
$ws.Range('E35').ClearFormats()
$ws.Range('D36').ClearFormats()
$ws.Range('E36').ClearFormats()
$ws.Range('E37').ClearFormats()
$ws.Range('D38').ClearFormats()
$ws.Range('E38').ClearFormats()
$ws.Range('E39').ClearFormats()
$ws.Range('E40').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('E41').ClearFormats()
$ws.Range('D42').ClearFormats()
$ws.Range('E42').ClearFormats()
$ws.Range('D43').ClearFormats()
$ws.Range('E43').ClearFormats()
$ws.Range('D44').ClearFormats()
$ws.Range('E44').ClearFormats()
$ws.Range('D45').ClearFormats()
$ws.Range('E45').ClearFormats()
$ws.Range('D46').ClearFormats()
$ws.Range('E46').ClearFormats()
$ws.Range('E47').ClearFormats()
$ws.Range('D48').ClearFormats()
$ws.Range('E48').ClearFormats()
$ws.Range('D49').ClearFormats()
$ws.Range('E49').ClearFormats()
$ws.Range('D50').ClearFormats()
$ws.Range('E50').ClearFormats()
$ws.Range('E51').ClearFormats()
